$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new values (previously held row 3's values, now updated further)
$ws.Range("B2").Value = 0.122138312895883
$ws.Range("C2").Value = 0.122138312895883
$ws.Range("D2").Value = 0.117239193324036
$ws.Range("E2").Value = 0.00077433319931779
$ws.Range("F2").Value = 0.7322

# Row 3: new values (previously held row 2's values, now slightly updated)
$ws.Range("B3").Value = 7.34768293483517
$ws.Range("C3").Value = 7.34768293483517
$ws.Range("D3").Value = 7.05295823772511
$ws.Range("E3").Value = 0.0465828837782763
$ws.Range("F3").Value = 0.0062

# Row 4: new values (same row, slightly updated precision)
$ws.Range("B4").Value = 0.246330117500345
$ws.Range("C4").Value = 0.246330117500345
$ws.Range("D4").Value = 0.236449510251333
$ws.Range("E4").Value = 0.00156168513752902
$ws.Range("F4").Value = 0.6317
